# Update database and change read_price algorithm
# - Drop the "1396/12" twelve-month period column and add a new
#   "1401/12" twelve-month period column (the whole five-year window
#   of period headers/data rolls forward by one year).
# - Refresh the figures for every expense/personnel row to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# ---------------------------------------------------------------
# Period headers (row 8 and row 24), columns E:I
# Old: 1396/12, 1397/12, 1398/12, 1399/12, 1400/12
# New: 1397/12, 1398/12, 1399/12, 1400/12, 1401/12
# ---------------------------------------------------------------
$years = @("1397", "1398", "1399", "1400", "1401")
$cols = @("E", "F", "G", "H", "I")

for ($i = 0; $i -lt 5; $i++) {
    $label = "دوازده ماهه منتهی به " + $years[$i] + "/12"
    $ws.Range($cols[$i] + "8").Value = $label
    $ws.Range($cols[$i] + "24").Value = $label
}

# ---------------------------------------------------------------
# Data rows: shift each figure one column to the left (E<-F, F<-G,
# G<-H, H<-I) and populate the new right-most column (I) with the
# newly reported 1401/12 figure.
# ---------------------------------------------------------------
$rowData = @{
    10 = @(14265, 23941, 141256, 7492, 62306)
    13 = @(5284, 2750, 7012, 10107, 1761)
    15 = @(0, 0, 0, 2598, 1286)
    16 = @(2238, 2204, 2470, 2711, 3177)
    17 = @(30646, 63066, 41735, 58066, 79324)
    19 = @(36103, 28961, 36441, 77027, 108352)
    20 = @(88536, 120922, 228914, 158001, 256206)
    26 = @(111, 95, 95, 105, 106)
    27 = @(294, 275, 256, 271, 318)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt 5; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}
